# Powerpoint writer: consolidate text run nodes.
# Merge the "First"/" " runs into a single "First " run on slide 1's title,
# and the "Third"/" " runs into a single "Third " run on slide 3's title,
# leaving the trailing "slide" run untouched.

$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(1, 6).Text = "First "

$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Characters(1, 6).Text = "Third "
